$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 52
$ws.Cells.Item($row, 1).Value = "AppyFizz20"
$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 20
$ws.Cells.Item($row, 4).Value = "AppyFizz20.jpg"

$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D53").Select()
